$d = $word.ActiveDocument

$pairs = @(
    @("209×4=", "701×5="),
    @("523×8=", "970×2="),
    @("737×2=", "151×9="),
    @("744×2=", "991×7="),
    @("270×2=", "172×4="),
    @("571×7=", "943×2="),
    @("976×7=", "592×6="),
    @("132×3=", "398×5="),
    @("245×9=", "637×7="),
    @("540×3=", "665×2="),
    @("695×8=", "817×7="),
    @("520×9=", "950×5="),
    @("565×8=", "908×4="),
    @("202×8=", "731×9="),
    @("440×4=", "432×6="),
    @("352×8=", "999×6="),
    @("336×2=", "616×7="),
    @("310×4=", "936×8="),
    @("135×5=", "870×8="),
    @("266×8=", "232×8="),
    @("681×9=", "423×4="),
    @("996×2=", "766×9="),
    @("830×5=", "231×6="),
    @("238×6=", "937×9="),
    @("479×5=", "489×7=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
